$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 16.602
$ws.Range("B3").Value = 5.484
$ws.Range("B14").Value = 5.755999999999999
$ws.Range("B16").Value = 5.63
$ws.Range("E18").Value = 16.481
$ws.Range("B21").Value = 8.996
$ws.Range("B23").Value = 7.849000000000001
$ws.Range("E24").Value = 16.547
$ws.Range("B25").Value = 5.67
$ws.Range("E25").Value = 17.164
$ws.Range("B26").Value = 6.397
$ws.Range("E27").Value = 16.659
$ws.Range("B29").Value = 5.392999999999999
$ws.Range("E30").Value = 16.466
$ws.Range("E31").Value = 16.346
$ws.Range("E39").Value = 16.218
$ws.Range("B40").Value = 9.137
$ws.Range("E42").Value = 16.663
$ws.Range("E48").Value = 17.39
$ws.Range("E51").Value = 16.617
$ws.Range("E52").Value = 16.543
$ws.Range("B53").Value = 5.89
$ws.Range("E55").Value = 16.416
$ws.Range("E56").Value = 16.214
$ws.Range("B57").Value = 5.091
$ws.Range("E57").Value = 16.553
$ws.Range("B59").Value = 4.435
$ws.Range("E60").Value = 16.569
$ws.Range("B65").Value = 5.737
$ws.Range("B69").Value = 5.106
$ws.Range("E73").Value = 16.735
$ws.Range("E74").Value = 16.607
$ws.Range("B79").Value = 5.548999999999999
$ws.Range("B83").Value = 5.523999999999999
$ws.Range("E89").Value = 17.374
$ws.Range("E90").Value = 16.348
$ws.Range("B91").Value = 6.318000000000001
$ws.Range("E92").Value = 17.7
$ws.Range("B93").Value = 5.635000000000001
$ws.Range("B100").Value = 5.728
